$wb = $excel.ActiveWorkbook

# New timestamp values for column F, rows 2-30 (both "Company" and "Driver 1" sheets)
$newTimestamps = @(
    45259.49596148148,
    45259.49619430555,
    45259.49622400463,
    45259.49657850694,
    45259.4968712037,
    45259.496924375,
    45259.497244756945,
    45259.49751269676,
    45259.497580972224,
    45259.4978718287,
    45259.49818141204,
    45259.498268773146,
    45259.49860806713,
    45259.49883331019,
    45259.49896023148,
    45259.49932715278,
    45259.49959344907,
    45259.499716238424,
    45259.500064872685,
    45259.500276631945,
    45259.50042986111,
    45259.50081430555,
    45259.50101581019,
    45259.50117337963,
    45259.5014859838,
    45259.50158930555,
    45259.501689398145,
    45259.50185825231,
    45259.50193409722,
)

$sheetNames = @("Company", "Driver 1")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 6).Value = $newTimestamps[$i]
    }
}
